# Apply odds updates to Sheet1 per commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$changes = @{
    "H2"  = 7
    "R2"  = 1.56
    "U2"  = 2.14
    "AE2" = 95
    "F3"  = 2.04
    "P3"  = 2.34
    "S3"  = 2.74
    "X3"  = 22
    "AD3" = 16
    "AJ3" = 25
    "AO3" = 50
    "G5"  = 1.63
    "H5"  = 4.6
    "P5"  = 1.16
    "N6"  = 3.55
    "P6"  = 1.81
    "AB6" = 9.800000000000001
    "AE6" = 1000
    "AF6" = 14.5
    "AI6" = 65
    "AO6" = 100
    "F7"  = 2.56
    "G7"  = 2.6
    "H7"  = 2.7
    "I7"  = 2.72
    "S7"  = 2.32
    "AA7" = 44
    "F8"  = 10
    "H8"  = 1.37
    "J8"  = 5.6
    "K8"  = 5.9
    "S8"  = 2.56
    "T8"  = 1.96
    "P9"  = 2.88
    "R9"  = 1.77
    "S9"  = 2.24
    "AL9" = 130
    "G10" = 1.3
    "I10" = 13
    "K10" = 7.2
    "N10" = 8.199999999999999
    "Q10" = 1.38
    "R10" = 2.02
    "S10" = 1.94
    "T10" = 1.75
    "U10" = 2.24
    "Z10" = 150
    "AA10" = 420
    "AC10" = 17.5
    "AD10" = 44
    "AO10" = 150
    "T11" = 1.9
    "AE11" = 17
    "AF11" = 55
    "AG11" = 23
    "AJ11" = 200
    "H12" = 2.46
    "P12" = 2.22
    "Q12" = 1.77
    "S12" = 2.92
    "U12" = 2.44
    "AF12" = 23
    "AN12" = 24
    "AO12" = 16.5
}

foreach ($addr in $changes.Keys) {
    $ws.Range($addr).Value = $changes[$addr]
}
